$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 220; this shifts old rows 220-273 down to 221-274
$ws.Rows.Item(220).Insert()

# Populate the newly inserted row 220 with the new data record
$ws.Cells.Item(220, 1).Value = 6
$ws.Cells.Item(220, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(220, 3).Value = "Metropolitana"
$ws.Cells.Item(220, 4).Value = 44889
$ws.Cells.Item(220, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(220, 5).Value = 13
$ws.Cells.Item(220, 6).Value = 100112022
$ws.Cells.Item(220, 7).Value = "Arveja Verde"
$ws.Cells.Item(220, 8).Value = "Sin especificar"
$ws.Cells.Item(220, 9).Value = "Primera"
$ws.Cells.Item(220, 10).Value = 290
$ws.Cells.Item(220, 11).Value = 20000
$ws.Cells.Item(220, 12).Value = 22000
$ws.Cells.Item(220, 13).Value = 21172
$ws.Cells.Item(220, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(220, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(220, 16).Value = 847
$ws.Cells.Item(220, 17).Value = 25
$ws.Cells.Item(220, 18).Value = "Hortaliza"
